$d = $word.ActiveDocument
$sec = $d.Sections(1)
$f1 = $sec.Footers(1)
$f1.Range.Find.Execute("QDAM", $true, $false, $false, $false, $false, $true, 1, $false, "QDAM2", 2)
Write-Output "done"
